# "Test.log to Home Page"
# Update the CreateSTP test-data sheet: the "Community Organizer" expert name
# moved from Kushalappa to ukumar1 (with display name "Uday Kumar"), and the
# mandatory/all-fields scenario test data was rebased from the *7011/*7012
# suffixed rows onto fresh *7016 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateSTP")

# Row 2 (Full Name) / Row 3 (Short Name): valid-mandatory (F) and valid-all (G)
# scenarios now point at the *7016 test rows instead of *7011/*7012.
$ws.Range("F2").Value = "FullName117016"
$ws.Range("F3").Value = "Short117016"
$ws.Range("G2").Value = "FullName217016"
$ws.Range("G3").Value = "Short217016"

# Row 5 (Community Organizer): replace the old expert "Kushalappa" with the
# "ukumar1" username already used elsewhere in the row.
$ws.Range("D5").Value = "ukumar1"
$ws.Range("E5").Value = "ukumar1"

# Update the saved view: selection moved to G3 and normal-view zoom is pinned
# at 100%.
$ws.Range("G3").Select()
$excel.ActiveWindow.Zoom = 100
